$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "description"

# Data rows
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Aluno"
$ws.Range("C2").Value = "Ambiente virtual para o aluno estudar."

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Professor"
$ws.Range("C3").Value = "Verificar a nota dos Alunos"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Diretor"
$ws.Range("C4").Value = "Verificar eficiências dos professores."

# Format the header cell A1 with bold font, thin box border, centered
# horizontally and top-aligned vertically, then propagate the same
# formatting to B1:C1 by copying A1's format only (keeps the style table
# minimal, matching one new font/border/cellXf).
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

$a1.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
